# Weekly update: insert new "Fruta" (Mango) price rows reported for
# Feria Lagunitas de Puerto Montt, interspersed among the existing
# chronologically-unsorted rows, plus one new row appended at the end.
#
# Net effect vs. before.xlsx:
#   - row 24 (new): 2021-05-07, vol 120, 10000/11000/10500, Peru, 2625
#   - row 39 (new): 2021-05-11, vol 120, 10000/11000/10500, Peru, 2625
#   - row 41 (new): 2021-05-14, vol 120, 10000/11000/10500, Peru, 2625
#   - row 65 (new, appended): 2021-05-18, vol 140, 10000/11000/10500, Peru, 2625
# All previously existing rows from row 24 onward shift down to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-MangoRow {
    param(
        [int]$RowNum,
        [string]$Fecha,
        [double]$Volumen,
        [double]$PrecioMin,
        [double]$PrecioMax,
        [double]$PrecioProm,
        [double]$PrecioKg
    )

    $ws.Cells.Item($RowNum, 1).Value = 4
    $ws.Cells.Item($RowNum, 2).Value = "Feria Lagunitas de Puerto Montt"
    $ws.Cells.Item($RowNum, 3).Value = "Los Lagos"
    $ws.Cells.Item($RowNum, 4).Value = $Fecha
    $ws.Cells.Item($RowNum, 5).Value = 10
    $ws.Cells.Item($RowNum, 6).Value = "Fruta"
    $ws.Cells.Item($RowNum, 7).Value = 100108
    $ws.Cells.Item($RowNum, 8).Value = "Tropicales y subtropicales"
    $ws.Cells.Item($RowNum, 9).Value = 100108002
    $ws.Cells.Item($RowNum, 10).Value = "Mango"
    $ws.Cells.Item($RowNum, 11).Value = "Sin especificar"
    $ws.Cells.Item($RowNum, 12).Value = "Primera"
    $ws.Cells.Item($RowNum, 13).Value = $Volumen
    $ws.Cells.Item($RowNum, 14).Value = $PrecioMin
    $ws.Cells.Item($RowNum, 15).Value = $PrecioMax
    $ws.Cells.Item($RowNum, 16).Value = $PrecioProm
    $ws.Cells.Item($RowNum, 17).Value = "`$/bandeja 4 kilos"
    $ws.Cells.Item($RowNum, 18).Value = "Perú"
    $ws.Cells.Item($RowNum, 19).Value = $PrecioKg
    $ws.Cells.Item($RowNum, 20).Value = 4
}

# 1) Insert a new row at 24 (everything from old row 24 downward shifts to 25+)
$ws.Rows.Item(24).Insert()
Set-MangoRow 24 "2021-05-07" 120 10000 11000 10500 2625

# 2) Insert a new row at 39 (old row 38 -- now sitting at row 39 after step 1 -- shifts to 40+)
$ws.Rows.Item(39).Insert()
Set-MangoRow 39 "2021-05-11" 120 10000 11000 10500 2625

# 3) Insert a new row at 41 (old row 39 -- now sitting at row 41 after steps 1-2 -- shifts to 42+)
$ws.Rows.Item(41).Insert()
Set-MangoRow 41 "2021-05-14" 120 10000 11000 10500 2625

# 4) Append one more new row at the very end (row 65)
Set-MangoRow 65 "2021-05-18" 140 10000 11000 10500 2625
# Row 65 is beyond the sheet's previous used range, so it doesn't inherit the
# "Fecha" column's date format automatically -- copy it from the column above.
$ws.Cells.Item(65, 4).NumberFormat = $ws.Cells.Item(64, 4).NumberFormat
